$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H74").Value = 4000
$ws.Range("I74").Value = 4000
$ws.Range("J74").Value = 0
$ws.Range("K74").Value = 4000
$ws.Range("L74").Value = 0
$ws.Range("M74").Value = -3064
$ws.Range("N74").ClearContents()
$ws.Range("H77").Value = 4000
$ws.Range("I77").Value = 4000
$ws.Range("J77").Value = 0
$ws.Range("K77").Value = 20000
$ws.Range("L77").Value = 0
$ws.Range("M77").Value = -15320
$ws.Range("N77").ClearContents()
$ws.Range("H86").Value = 60054.35
$ws.Range("I86").Value = 67888.266
$ws.Range("J86").Value = 1300
$ws.Range("K86").Value = 67888.266
$ws.Range("L86").Value = 1300
$ws.Range("M86").Value = -66765.266
$ws.Range("N86").Value = -3546
$ws.Range("H89").Value = 60054.35
$ws.Range("I89").Value = 67888.266
$ws.Range("J89").Value = 1300
$ws.Range("K89").Value = 339441.33
$ws.Range("L89").Value = 6500
$ws.Range("M89").Value = -333825.33
$ws.Range("N89").Value = -17732
$ws.Range("H98").Value = 2280.4
$ws.Range("I98").Value = 2157.5715
$ws.Range("K98").Value = 2157.5715
$ws.Range("M98").Value = -659.5715
$ws.Range("H113").Value = 1925
$ws.Range("I113").Value = 1925
$ws.Range("J113").Value = 0
$ws.Range("K113").Value = 1925
$ws.Range("L113").Value = 0
$ws.Range("M113").Value = 1329
$ws.Range("N113").ClearContents()
$ws.Range("H116").Value = 2413.125
$ws.Range("I116").Value = 1505
$ws.Range("J116").Value = 2542.8572
$ws.Range("K116").Value = 1505
$ws.Range("L116").Value = 2542.8572
$ws.Range("M116").Value = 1937
$ws.Range("N116").Value = -9426.8572
$ws.Range("H122").Value = 2280.4
$ws.Range("I122").Value = 2157.5715
$ws.Range("K122").Value = 6472.7145
$ws.Range("M122").Value = -4022.7145
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H5").Value = 0
$ws.Range("I5").Value = 0
$ws.Range("K5").Value = 0
$ws.Range("M5").ClearContents()
$ws.Range("H25").Value = 7503.2
$ws.Range("I25").Value = 758
$ws.Range("K25").Value = 758
$ws.Range("M25").Value = -356
$ws.Range("H32").Value = 14200.337
$ws.Range("I32").Value = 14446.695
$ws.Range("K32").Value = 14446.695
$ws.Range("M32").Value = -14159.695
$ws.Range("H61").Value = 2531.8096
$ws.Range("I61").Value = 1982.3334
$ws.Range("K61").Value = 1982.3334
$ws.Range("M61").Value = -1770.3334
$ws.Range("H63").Value = 911600.4399999999
$ws.Range("I63").Value = 1252225.6
$ws.Range("J63").Value = 3266.6667
$ws.Range("K63").Value = 1252225.6
$ws.Range("L63").Value = 3266.6667
$ws.Range("M63").Value = -1251539.6
$ws.Range("N63").Value = -4638.6667
$ws.Range("H66").Value = 911600.4399999999
$ws.Range("I66").Value = 1252225.6
$ws.Range("J66").Value = 3266.6667
$ws.Range("K66").Value = 6261128
$ws.Range("L66").Value = 16333.3335
$ws.Range("M66").Value = -6257696
$ws.Range("N66").Value = -23197.3335
$ws.Range("H132").Value = 4358.2617
$ws.Range("I132").Value = 5641.4546
$ws.Range("J132").Value = 2946.75
$ws.Range("K132").Value = 16924.3638
$ws.Range("L132").Value = 8840.25
$ws.Range("M132").Value = -14394.3638
$ws.Range("N132").Value = -13900.25
$ws.Range("H136").Value = 2531.8096
$ws.Range("I136").Value = 1982.3334
$ws.Range("K136").Value = 5947.0002
$ws.Range("M136").Value = -3397.0002
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 0
$ws.Range("I4").Value = 0
$ws.Range("K4").Value = 0
$ws.Range("M4").ClearContents()
$ws.Range("H20").Value = 1791.6086
$ws.Range("I20").Value = 1883.3529
$ws.Range("J20").Value = 1531.6666
$ws.Range("K20").Value = 1883.3529
$ws.Range("L20").Value = 1531.6666
$ws.Range("M20").Value = -1636.3529
$ws.Range("N20").Value = -2025.6666
$ws.Range("H22").Value = 411.875
$ws.Range("I22").Value = 365.26666
$ws.Range("K22").Value = 365.26666
$ws.Range("M22").Value = -192.26666
$ws.Range("H76").Value = 0
$ws.Range("J76").Value = 0
$ws.Range("L76").Value = 0
$ws.Range("N76").ClearContents()
$ws.Range("H79").Value = 0
$ws.Range("J79").Value = 0
$ws.Range("L79").Value = 0
$ws.Range("N79").ClearContents()
$ws.Range("H86").Value = 1983.8334
$ws.Range("I86").Value = 1358
$ws.Range("J86").Value = 2860
$ws.Range("K86").Value = 1358
$ws.Range("L86").Value = 2860
$ws.Range("M86").Value = -235
$ws.Range("N86").Value = -5106
$ws.Range("H89").Value = 1983.8334
$ws.Range("I89").Value = 1358
$ws.Range("J89").Value = 2860
$ws.Range("K89").Value = 6790
$ws.Range("L89").Value = 14300
$ws.Range("M89").Value = -1174
$ws.Range("N89").Value = -25532
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 5209809.5
$ws.Range("I31").Value = 2108.5833
$ws.Range("K31").Value = 2108.5833
$ws.Range("M31").Value = -1813.5833
$ws.Range("H34").Value = 5209809.5
$ws.Range("I34").Value = 2108.5833
$ws.Range("K34").Value = 2108.5833
$ws.Range("M34").Value = -1906.5833
$ws.Range("H58").Value = 907.3570999999999
$ws.Range("I58").Value = 931.2381
$ws.Range("J58").Value = 835.7143
$ws.Range("K58").Value = 931.2381
$ws.Range("L58").Value = 835.7143
$ws.Range("M58").Value = -728.2381
$ws.Range("N58").Value = -1241.7143
$ws.Range("H99").Value = 2107.0557
$ws.Range("I99").Value = 1576.5
$ws.Range("J99").Value = 2531.5
$ws.Range("K99").Value = 1576.5
$ws.Range("L99").Value = 2531.5
$ws.Range("M99").Value = -78.5
$ws.Range("N99").Value = -5527.5
$ws.Range("H126").Value = 2107.0557
$ws.Range("I126").Value = 1576.5
$ws.Range("J126").Value = 2531.5
$ws.Range("K126").Value = 4729.5
$ws.Range("L126").Value = 7594.5
$ws.Range("M126").Value = -2259.5
$ws.Range("N126").Value = -12534.5
$ws.Range("H132").Value = 3630.4443
$ws.Range("I132").Value = 3406.6667
$ws.Range("J132").Value = 3854.2222
$ws.Range("K132").Value = 10220.0001
$ws.Range("L132").Value = 11562.6666
$ws.Range("M132").Value = -7690.000100000001
$ws.Range("N132").Value = -16622.6666
$ws.Range("H134").Value = 904.3913
$ws.Range("I134").Value = 915.05
$ws.Range("K134").Value = 2745.15
$ws.Range("M134").Value = -210.1499999999996
$ws.Range("H136").Value = 907.3570999999999
$ws.Range("I136").Value = 931.2381
$ws.Range("J136").Value = 835.7143
$ws.Range("K136").Value = 2793.7143
$ws.Range("L136").Value = 2507.1429
$ws.Range("M136").Value = -243.7143000000001
$ws.Range("N136").Value = -7607.1429
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H97").Value = 288.57144
$ws.Range("I97").Value = 288.57144
$ws.Range("J97").Value = 0
$ws.Range("K97").Value = 865.71432
$ws.Range("L97").Value = 0
$ws.Range("M97").Value = -369.71432
$ws.Range("N97").ClearContents()
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 23187360
$ws.Range("I70").Value = 28337894
$ws.Range("J70").Value = 9950
$ws.Range("K70").Value = 28337894
$ws.Range("L70").Value = 9950
$ws.Range("M70").Value = -28337624
$ws.Range("N70").Value = -10490
$ws.Range("H73").Value = 23187360
$ws.Range("I73").Value = 28337894
$ws.Range("J73").Value = 9950
$ws.Range("K73").Value = 28337894
$ws.Range("L73").Value = 9950
$ws.Range("M73").Value = -28336958
$ws.Range("N73").Value = -11822
$ws.Range("H80").Value = 5407.04
$ws.Range("I80").Value = 4759.8667
$ws.Range("J80").Value = 6377.8
$ws.Range("K80").Value = 4759.8667
$ws.Range("L80").Value = 6377.8
$ws.Range("M80").Value = -3761.8667
$ws.Range("N80").Value = -8373.799999999999
$ws.Range("H83").Value = 5407.04
$ws.Range("I83").Value = 4759.8667
$ws.Range("J83").Value = 6377.8
$ws.Range("K83").Value = 23799.3335
$ws.Range("L83").Value = 31889
$ws.Range("M83").Value = -18807.3335
$ws.Range("N83").Value = -41873
$ws.Range("H102").Value = 1326.6957
$ws.Range("I102").Value = 1280.091
$ws.Range("J102").Value = 1369.4166
$ws.Range("K102").Value = 1280.091
$ws.Range("L102").Value = 1369.4166
$ws.Range("M102").Value = 341.9090000000001
$ws.Range("N102").Value = -4613.4166
$ws.Range("H122").Value = 6992.3076
$ws.Range("I122").Value = 29050
$ws.Range("K122").Value = 87150
$ws.Range("M122").Value = -84700
$ws.Range("H132").Value = 34023.97
$ws.Range("I132").Value = 40721.883
$ws.Range("J132").Value = 4999.6665
$ws.Range("K132").Value = 122165.649
$ws.Range("L132").Value = 14998.9995
$ws.Range("M132").Value = -119635.649
$ws.Range("N132").Value = -20058.9995
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 2611.5557
$ws.Range("I40").Value = 2250.6667
$ws.Range("J40").Value = 3333.3333
$ws.Range("K40").Value = 2250.6667
$ws.Range("L40").Value = 3333.3333
$ws.Range("M40").Value = -2114.6667
$ws.Range("N40").Value = -3605.3333
$ws.Range("H46").Value = 2089.7917
$ws.Range("I46").Value = 1955.6666
$ws.Range("J46").Value = 2170.2666
$ws.Range("K46").Value = 1955.6666
$ws.Range("L46").Value = 2170.2666
$ws.Range("M46").Value = -1767.6666
$ws.Range("N46").Value = -2546.2666
$ws.Range("H122").Value = 5800.8
$ws.Range("J122").Value = 3000
$ws.Range("L122").Value = 9000
$ws.Range("N122").Value = -13900
$ws.Range("H132").Value = 11202.4
$ws.Range("I132").Value = 16243.357
$ws.Range("J132").Value = 4786.636
$ws.Range("K132").Value = 48730.071
$ws.Range("L132").Value = 14359.908
$ws.Range("M132").Value = -46200.071
$ws.Range("N132").Value = -19419.908
$ws.Range("H139").Value = 47800
$ws.Range("I139").Value = 0
$ws.Range("J139").Value = 47800
$ws.Range("K139").Value = 0
$ws.Range("L139").Value = 47800
$ws.Range("M139").ClearContents()
$ws.Range("N139").Value = -58080
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H63").Value = 10000
$ws.Range("J63").Value = 10000
$ws.Range("L63").Value = 10000
$ws.Range("N63").Value = -11248
$ws.Range("H66").Value = 10000
$ws.Range("J66").Value = 10000
$ws.Range("L66").Value = 30000
$ws.Range("N66").Value = -36240
